$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (8 and 9) before the current row 8 ("extr1"),
# pushing the existing extr1..extr8 rows down by two rows.
$ws.Rows.Item(8).Resize(2, 1).EntireRow.Insert()

# Re-apply the index-column formatting (A7 -> A8:A9) that the row
# insert otherwise fails to carry down correctly.
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

# New row 8: line7
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# New row 9: line8
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $false

# Renumber the index column (A) and the in_service flags (E) for the
# shifted extr1..extr8 rows, now located at rows 10..17.
for ($i = 0; $i -lt 8; $i++) {
    $row = 10 + $i
    $ws.Cells.Item($row, 1).Value = 8 + $i
}

# extr1 (row 10) and extr2 (row 11) now become in service.
$ws.Cells.Item(10, 5).Value = $true
$ws.Cells.Item(11, 5).Value = $true
